$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new rows at the right spots, shifting existing rows down ---
# Insert a row for "Mazda MX 30 2020 " before current row 3 (Honda Jazz)
$ws.Rows.Item(3).Insert()
# Insert a row for "Landrover Defender " before current row 5 (SEAT Leon)
$ws.Rows.Item(5).Insert()
# Insert two rows for "Honda e " and "Hyundai i10 2020 " before current row 8 (Isuzu D-Max)
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# --- Widen column A from 20 to 21 characters ---
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668

# --- Row 2: Toyota Yaris (text reformatted, values unchanged) ---
$ws.Range("A2").Value = "Toyota Yaris 2020 "

# --- Row 3: Mazda MX 30 2020 (new row) ---
$ws.Range("A3").Value = "Mazda MX 30 2020 "
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 46.99
$ws.Range("D3").Value = 36.75
$ws.Range("E3").Value = 5.42
$ws.Range("F3").Value = 7.23
$ws.Range("G3").Value = 3.61
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 100

# --- Row 4: Honda Jazz (text reformatted, values unchanged) ---
$ws.Range("A4").Value = "Honda Jazz 2020 "

# --- Row 5: Landrover Defender (new row) ---
$ws.Range("A5").Value = "Landrover Defender "
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 49.68
$ws.Range("D5").Value = 17.2
$ws.Range("E5").Value = 22.29
$ws.Range("F5").Value = 6.37
$ws.Range("G5").Value = 4.46
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 100

# --- Row 6: SEAT Leon (text reformatted, values unchanged) ---
$ws.Range("A6").Value = "SEAT Leon 2020 "

# --- Row 7: Kia Sorento (text reformatted, values unchanged) ---
$ws.Range("A7").Value = "Kia Sorento 2020 "

# --- Row 8: Honda e (new row) ---
$ws.Range("A8").Value = "Honda e "
$ws.Range("B8").Value = 15.13
$ws.Range("C8").Value = 36.84
$ws.Range("D8").Value = 19.74
$ws.Range("E8").Value = 9.869999999999999
$ws.Range("F8").Value = 3.95
$ws.Range("G8").Value = 9.21
$ws.Range("H8").Value = 5.26
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 100

# --- Row 9: Hyundai i10 2020 (new row) ---
$ws.Range("A9").Value = "Hyundai i10 2020 "
$ws.Range("B9").Value = 20.99
$ws.Range("C9").Value = 8.02
$ws.Range("D9").Value = 20.37
$ws.Range("E9").Value = 22.22
$ws.Range("F9").Value = 8.64
$ws.Range("G9").Value = 12.35
$ws.Range("H9").Value = 7.41
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 100

# --- Row 10: Isuzu D Max (text reformatted, values unchanged) ---
$ws.Range("A10").Value = "Isuzu D Max 2020 "

# --- Row 11: Audi A3 (new row, appended) ---
$ws.Range("A11").Value = "Audi A3 "
$ws.Range("B11").Value = 1.81
$ws.Range("C11").Value = 22.29
$ws.Range("D11").Value = 25.3
$ws.Range("E11").Value = 21.69
$ws.Range("F11").Value = 12.65
$ws.Range("G11").Value = 12.65
$ws.Range("H11").Value = 3.61
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 100
